$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the selected cell on the sheet view
$ws.Range("D4").Select()

# Update the "ligne" (line range) column B values for rows 7-17
# (order matters for shared-string table layout, matching the authored diff)
$ws.Range("B9").Value = "43 à 75"
$ws.Range("B10").Value = "81 à 112"
$ws.Range("B11").Value = "117 à 123"
$ws.Range("B12").Value = "126 à 186"
$ws.Range("B13").Value = "189 à 265"
$ws.Range("B14").Value = "280 à 398"
$ws.Range("B15").Value = "408 à 414"
$ws.Range("B16").Value = "424 à 430"
$ws.Range("B17").Value = "440 à 447"
$ws.Range("B7").Value = "67 à 137"
$ws.Range("B8").Value = "140 à 151"
